$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$theme = $m.Theme
try {
  $cs = $theme.ThemeColorScheme
  Write-Host "ThemeColorScheme: $cs"
} catch { Write-Host "ERR1: $_" }
